$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Columns("A").ColumnWidth = 17.6640625
$ws.Columns("A").BestFit = $true
